$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Collect the cells that must keep their original text/string format (not auto-converted to numbers)
$updates = @{
    'D2' = '279.25'
    'E2' = '6.11%'
    'D3' = '27.03'
    'E3' = '1.37%'
    'D4' = '4.924'
    'E4' = '4.96%'
    'D5' = '0.06362'
    'E5' = '4.33%'
    'D6' = '6.943'
    'E6' = '3.50%'
    'D7' = '3.360'
    'E7' = '6.15%'
    'D8' = '0.8846'
    'E8' = '4.02%'
    'D9' = '0.9464'
    'E9' = '3.84%'
    'D10' = '0.1470'
    'E10' = '4.31%'
    'D11' = '0.05117'
    'E11' = '7.29%'
    'D12' = '0.07406'
    'E12' = '4.47%'
    'D13' = '0.03156'
    'E13' = '0.71%'
    'D14' = '0.09040'
    'E14' = '-0.09%'
    'D15' = '0.001559'
    'E15' = '2.02%'
    'D16' = '0.0006272'
    'E16' = '1.31%'
    'D17' = '0.005809'
    'D18' = '3.486'
    'D19' = '2.297'
    'E19' = '7.03%'
    'D20' = '0.3098'
    'E20' = '0.84%'
    'D21' = '0.1302'
    'E21' = '1.71%'
    'D22' = '3.885'
    'E22' = '-5.50%'
    'D23' = '0.04340'
    'E23' = '2.64%'
    'D24' = '0.001176'
    'E24' = '-0.08%'
    'D25' = '0.003638'
    'D26' = '0.0001200'
    'E26' = '0.03%'
    'D27' = '0.0001694'
    'E27' = '-12.55%'
    'D40' = '0.04065'
    'E40' = '3.46%'
    'D41' = '0.006627'
    'E41' = '58.84%'
    'D42' = '0.1167'
    'E42' = '4.77%'
    'D43' = '0.002350'
    'E43' = '11.43%'
    'D44' = '0.01253'
    'E44' = '8.48%'
    'D45' = '0.00005266'
    'E45' = '3.64%'
    'E46' = '0.03%'
    'D47' = '2.374'
    'E47' = '821.51%'
    'D48' = '0.02261'
    'E48' = '6.56%'
    'E49' = '0.03%'
    'E50' = '-0.04%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"   # force Text format so numeric-looking strings are not coerced to numbers/dates
    $cell.Value = $updates[$addr]
}
